# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" worksheet (fund-holding detail) positioned
#    right before the "总计" (summary) sheet.
# 2) Insert a new first data row into "总计" summarizing the 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the "2022-Q1" sheet, positioned after "2021-Q4" (i.e.
# right before "总计").
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $afterSheet)
$q1.Name = "2022-Q1"

# Borrow the row/column formatting (header style + index-column style)
# from the existing "2021-Q1" sheet, which has the identical 8-row,
# 8-column (A:H) layout we need here.
$template = $wb.Worksheets.Item("2021-Q1")
$template.Range("A1:H8").Copy($q1.Range("A1:H8"))

# Header row
$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

# Columns B-G hold text (fund codes / names / numeric-looking text such
# as "5.86"); force them to Text format first so codes like "002434"
# keep their leading zeros and figures stay literal text, matching the
# source data.
$q1.Range("B2:G8").NumberFormat = "@"

$q1Data = @(
    @(0, "002434", "中银宏利灵活配置混合A", "5.86", "30.78", "0.68", "0.0398", 7),
    @(1, "003966", "中银润利灵活配置混合A", "5.69", "25.96", "0.54", "0.0307", 8),
    @(2, "002261", "中银宝利灵活配置混合A", "4.35", "31.13", "0.67", "0.0291", 7),
    @(3, "160639", "鹏华中证高铁产业指数（LOF）", "0.89", "94.72", "2.66", "0.0237", 7),
    @(4, "003967", "中银润利灵活配置混合C", "3.85", "25.96", "0.54", "0.0208", 8),
    @(5, "002435", "中银宏利灵活配置混合C", "2.33", "30.78", "0.68", "0.0158", 7),
    @(6, "002262", "中银宝利灵活配置混合C", "1.60", "31.13", "0.67", "0.0107", 7)
)

$r = 2
foreach ($row in $q1Data) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).Value = $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: prepend a "2022-Q1" row to the "总计" summary sheet.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push the existing data rows down one row, then restore the
# second-row formatting (index-column style, etc.) onto the freshly
# inserted, now-blank row by copying it from the row beneath it.
$total.Rows.Item(2).Insert()
$total.Range("A3:D3").Copy($total.Range("A2:D2"))

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 7
$total.Cells.Item(2, 4).Value = 0.17

# The "A" column is a running 0-based index; bump every pre-existing
# row's index by one now that they have shifted down a row.
for ($row = 3; $row -le 7; $row++) {
    $total.Cells.Item($row, 1).Value = $row - 2
}
